# Donorset fixed, donorset population was added.
#
# The "top donor" (D129:D134) and "top party" (D136:D141) name lists in
# column D carried a stray trailing space on every entry; strip it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 129,130,131,132,133,134,136,137,138,139,140,141
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Value2 = $cell.Value2.TrimEnd()
}
